# Generate Report for handoff
#
# b.md.md has been handed off again for localization: its status moves
# from "Handed back: in sync with en-US" to "Ready for handoff", and a
# new handoff file + handoff datetime are recorded for both the zh-cn and
# de-de locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is b.md.md ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 is b.md.md ---
$zhcnFile = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcnUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4d182e58bbab71f159159a5b94ad54e8d284a83/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"

$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = $zhcnFile
$zhcn.Range("D3").Value = "2016-01-25 13:38:17"

$zhcnLink = $zhcn.Hyperlinks.Item(6)
$zhcnLink.TextToDisplay = $zhcnFile
$zhcnLink.Address = $zhcnUrl

# --- de-de sheet: row 3 is b.md.md ---
$dedeFile = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dedeUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e9653abe0271fcdf25080e1623c29f7ad7a3a56/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"

$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = $dedeFile
$dede.Range("D3").Value = "2016-01-25 13:38:26"

$dedeLink = $dede.Hyperlinks.Item(6)
$dedeLink.TextToDisplay = $dedeFile
$dedeLink.Address = $dedeUrl
